$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. MCD sheet: add "user_corpo" attribute to the `user` entity (G10)
#    and "task_type" attribute to the `task` entity (B14).
# ---------------------------------------------------------------------
$mcd = $wb.Worksheets.Item("MCD")
$mcd.Range("G10").Value = "user_corpo"
$mcd.Range("B14").Value = "task_type"

# Grow the two underlying Excel tables so the new attribute rows are
# included (Table1 = `user`, Table7 = `task`).
$userTable = $mcd.ListObjects.Item("Table1")
$userTable.Resize($mcd.Range("G7:G10"))

$taskTable = $mcd.ListObjects.Item("Table7")
$taskTable.Resize($mcd.Range("B7:B14"))

# Widen columns I and K so the (now longer) entity/attribute lists stay
# readable.
$mcd.Columns.Item(9).ColumnWidth = 36.333333333333336
$mcd.Columns.Item(11).ColumnWidth = 22.666666666666668

# ---------------------------------------------------------------------
# 2. Insert the new "task_type" lookup sheet right after "MCD".
# ---------------------------------------------------------------------
$taskType = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $mcd)
$taskType.Name = "task_type"

$taskType.Range("A1").Value = 1
$taskType.Range("B1").Value = "WS_BS"
$taskType.Range("A2").Value = 2
$taskType.Range("B2").Value = "WS_TS"
$taskType.Range("A3").Value = 3
$taskType.Range("B3").Value = "WS_FS"
$taskType.Range("A4").Value = 4
$taskType.Range("B4").Value = "URGENT"
$taskType.Range("A5").Value = 5
$taskType.Range("B5").Value = "DIVERS"

# The last two rows (URGENT / DIVERS) are shown in italics.
$taskType.Range("A4:B5").Font.Italic = $true

# Match the page setup used throughout the rest of the workbook.
$taskType.PageSetup.PaperSize = 9
$taskType.PageSetup.Orientation = 1

# MCD keeps its scroll position but the selection moved to B18.
$mcd.Range("B18").Select()

$taskType.Range("E7").Select()
$taskType.Activate()
